$wb = $excel.ActiveWorkbook

# --- Update the daily conversion text on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = $wsHoja1.Range("A1").Value()
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 1.78 = 6518.26 pesos"), "✅ 1000 Bs = 1.75 = 6386.74 pesos"
$text = $text -replace [regex]::Escape("✅ 6518.26 pesos = 1.77 = 945.96 Bs"), "✅ 6386.74 pesos = 1.74 = 930.53 Bs"
$wsHoja1.Range("A1").Value = $text

# --- Update the tasas sheet with new rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 570
$wsTasas.Range("O10").Value = 3640.44
$wsTasas.Range("N12").Value = 3672
$wsTasas.Range("O12").Value = 535
